$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"2.220712"
$ws.Range("H2").Value = [double]"6.662135999999999"
$ws.Range("I2").Value = [double]"0.004164179109543329"
$ws.Range("J2").Value = [double]"0.00416417910954333"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"0.6120643333333334"
$ws.Range("N2").Value = [double]"1.836193"
$ws.Range("O2").Value = [double]"0.01381035804569015"
$ws.Range("P2").Value = [double]"0.01381035804569015"
$ws.Range("Q2").Value = [double]"1.359218609805333"
$ws.Range("R2").Value = [double]"12.232967488248"
$ws.Range("S2").Value = [double]"5.750880446917658E-05"
$ws.Range("T2").Value = [double]"5.750880446917659E-05"
$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"2.220712"
$ws.Range("H3").Value = [double]"6.662135999999999"
$ws.Range("I3").Value = [double]"0.004164179109543329"
$ws.Range("J3").Value = [double]"0.00416417910954333"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"4.841148333333333"
$ws.Range("N3").Value = [double]"14.523445"
$ws.Range("O3").Value = [double]"0.1092336020815287"
$ws.Range("P3").Value = [double]"0.1092336020815287"
$ws.Range("Q3").Value = [double]"10.75079619761333"
$ws.Range("R3").Value = [double]"96.75716577851999"
$ws.Range("S3").Value = [double]"0.0004548682838480705"
$ws.Range("T3").Value = [double]"0.0004548682838480705"
$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"2.220712"
$ws.Range("H4").Value = [double]"6.662135999999999"
$ws.Range("I4").Value = [double]"0.004164179109543329"
$ws.Range("J4").Value = [double]"0.00416417910954333"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"38.86601"
$ws.Range("N4").Value = [double]"116.59803"
$ws.Range("O4").Value = [double]"0.8769560398727811"
$ws.Range("P4").Value = [double]"0.8769560398727813"
$ws.Range("Q4").Value = [double]"86.31021479911998"
$ws.Range("R4").Value = [double]"776.7919331920799"
$ws.Range("S4").Value = [double]"0.003651802021226082"
$ws.Range("T4").Value = [double]"0.003651802021226083"
$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"527.8012189999999"
$ws.Range("H5").Value = [double]"1583.403657"
$ws.Range("I5").Value = [double]"0.9897090708526379"
$ws.Range("J5").Value = [double]"0.9897090708526382"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"0.6120643333333334"
$ws.Range("N5").Value = [double]"1.836193"
$ws.Range("O5").Value = [double]"0.01381035804569015"
$ws.Range("P5").Value = [double]"0.01381035804569015"
$ws.Range("Q5").Value = [double]"323.0483012397557"
$ws.Range("R5").Value = [double]"2907.434711157801"
$ws.Range("S5").Value = [double]"0.01366823662954226"
$ws.Range("T5").Value = [double]"0.01366823662954226"
$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"527.8012189999999"
$ws.Range("H6").Value = [double]"1583.403657"
$ws.Range("I6").Value = [double]"0.9897090708526379"
$ws.Range("J6").Value = [double]"0.9897090708526382"
$ws.Range("K6").Value = [double]"3"
$ws.Range("M6").Value = [double]"4.841148333333333"
$ws.Range("N6").Value = [double]"14.523445"
$ws.Range("O6").Value = [double]"0.1092336020815287"
$ws.Range("P6").Value = [double]"0.1092336020815287"
$ws.Range("Q6").Value = [double]"2555.163991693151"
$ws.Range("R6").Value = [double]"22996.47592523836"
$ws.Range("S6").Value = [double]"0.1081094868219965"
$ws.Range("T6").Value = [double]"0.1081094868219966"
$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"527.8012189999999"
$ws.Range("H7").Value = [double]"1583.403657"
$ws.Range("I7").Value = [double]"0.9897090708526379"
$ws.Range("J7").Value = [double]"0.9897090708526382"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"38.86601"
$ws.Range("N7").Value = [double]"116.59803"
$ws.Range("O7").Value = [double]"0.8769560398727811"
$ws.Range("P7").Value = [double]"0.8769560398727813"
$ws.Range("Q7").Value = [double]"20513.52745566619"
$ws.Range("R7").Value = [double]"184621.7471009957"
$ws.Range("S7").Value = [double]"0.8679313474010991"
$ws.Range("T7").Value = [double]"0.8679313474010995"
$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"3.26733"
$ws.Range("H8").Value = [double]"9.80199"
$ws.Range("I8").Value = [double]"0.006126750037818593"
$ws.Range("J8").Value = [double]"0.006126750037818595"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"0.6120643333333334"
$ws.Range("N8").Value = [double]"1.836193"
$ws.Range("O8").Value = [double]"0.01381035804569015"
$ws.Range("P8").Value = [double]"0.01381035804569015"
$ws.Range("Q8").Value = [double]"1.99981615823"
$ws.Range("R8").Value = [double]"17.99834542407"
$ws.Range("S8").Value = [double]"8.461261167872047E-05"
$ws.Range("T8").Value = [double]"8.461261167872049E-05"
$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"3.26733"
$ws.Range("H9").Value = [double]"9.80199"
$ws.Range("I9").Value = [double]"0.006126750037818593"
$ws.Range("J9").Value = [double]"0.006126750037818595"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"4.841148333333333"
$ws.Range("N9").Value = [double]"14.523445"
$ws.Range("O9").Value = [double]"0.1092336020815287"
$ws.Range("P9").Value = [double]"0.1092336020815287"
$ws.Range("Q9").Value = [double]"15.81762918395"
$ws.Range("R9").Value = [double]"142.35866265555"
$ws.Range("S9").Value = [double]"0.0004548682838480705"
$ws.Range("T9").Value = [double]"0.0006692469756840673"
$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"3.26733"
$ws.Range("H10").Value = [double]"9.80199"
$ws.Range("I10").Value = [double]"0.006126750037818593"
$ws.Range("J10").Value = [double]"0.006126750037818595"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"38.86601"
$ws.Range("N10").Value = [double]"116.59803"
$ws.Range("O10").Value = [double]"0.8769560398727811"
$ws.Range("P10").Value = [double]"0.8769560398727813"
$ws.Range("Q10").Value = [double]"126.9880804533"
$ws.Range("R10").Value = [double]"1142.8927240797"
$ws.Range("S10").Value = [double]"0.005372890450455805"
$ws.Range("T10").Value = [double]"0.005372890450455808"
